# Updates cryptos list cell values to reflect the latest scrape (GitHub Actions run).
# All D (Price) and E (Volume/1h) columns are stored as text in this workbook, and some
# D values look numeric (e.g. "0.0950", "340.40"); we force those to stay text by using a
# leading apostrophe so Excel doesn't silently coerce them into numbers and strip formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'59.209.08"
$ws.Range("E2").Value = '  +1.61%  '

# Row 3
$ws.Range("D3").Value = "'2.586.85"
$ws.Range("E3").Value = '  -0.30%  '

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = "'522.43"
$ws.Range("E5").Value = '  -0.04%  '

# Row 6
$ws.Range("D6").Value = "'139.25"
$ws.Range("E6").Value = '  -3.66%  '

# Row 7
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("E8").Value = '  -0.40%  '

# Row 9
$ws.Range("D9").Value = "'2.598.93"
$ws.Range("E9").Value = '  -0.64%  '

# Row 10
$ws.Range("D10").Value = "'6.54"
$ws.Range("E10").Value = '  -1.78%  '

# Row 11
$ws.Range("E11").Value = '  -0.56%  '

# Row 12
$ws.Range("D12").Value = "'0.329"
$ws.Range("E12").Value = '  -1.68%  '

# Row 13
$ws.Range("E13").Value = '  +2.94%  '

# Row 14
$ws.Range("D14").Value = "'3.048.65"
$ws.Range("E14").Value = '  -0.21%  '

# Row 15
$ws.Range("D15").Value = "'58.975.01"
$ws.Range("E15").Value = '  +1.22%  '

# Row 16
$ws.Range("D16").Value = "'20.57"
$ws.Range("E16").Value = '  -0.03%  '

# Row 17
$ws.Range("D17").Value = "'2.593.64"
$ws.Range("E17").Value = '  -1.62%  '

# Row 18
$ws.Range("E18").Value = '  -1.23%  '

# Row 19
$ws.Range("D19").Value = "'340.40"
$ws.Range("E19").Value = '  +0.00%  '

# Row 20
$ws.Range("E20").Value = '  -1.58%  '

# Row 21
$ws.Range("E21").Value = '  -2.54%  '

# Row 22
$ws.Range("D22").Value = "'6.44"
$ws.Range("E22").Value = '  +0.48%  '

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  +0.15%  '

# Row 24
$ws.Range("E24").Value = '  +0.41%  '

# Row 25
$ws.Range("E25").Value = '  +0.89%  '

# Row 26
$ws.Range("E26").Value = '  -0.12%  '

# Row 27
$ws.Range("E27").Value = '  +0.20%  '

# Row 28
$ws.Range("E28").Value = '  -0.10%  '

# Row 30
$ws.Range("E30").Value = '  -4.00%  '

# Row 31
$ws.Range("D31").Value = "'5.88"
$ws.Range("E31").Value = '  -6.24%  '

# Row 32
$ws.Range("E32").Value = '  -0.60%  '

# Row 33
$ws.Range("D33").Value = "'18.69"
$ws.Range("E33").Value = '  -1.08%  '

# Row 34
$ws.Range("D34").Value = "'149.36"
$ws.Range("E34").Value = '  -0.28%  '

# Row 35
$ws.Range("D35").Value = "'3.96"
$ws.Range("E35").Value = '  -2.35%  '

# Row 36
$ws.Range("E36").Value = '  -2.68%  '

# Row 37
$ws.Range("D37").Value = "'36.73"
$ws.Range("E37").Value = '  +1.70%  '

# Row 38
$ws.Range("E38").Value = '  +0.24%  '

# Row 39
$ws.Range("D39").Value = "'0.823"
$ws.Range("E39").Value = '  -3.09%  '

# Row 40
$ws.Range("E40").Value = '  -6.48%  '

# Row 41
$ws.Range("E41").Value = '  -1.25%  '

# Row 42
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  +0.30%  '

# Row 43
$ws.Range("D43").Value = "'272.03"
$ws.Range("E43").Value = '  -0.75%  '

# Row 44
$ws.Range("D44").Value = "'0.602"
$ws.Range("E44").Value = '  +0.64%  '

# Row 45
$ws.Range("D45").Value = "'10.78"
$ws.Range("E45").Value = '  +1.40%  '

# Row 46
$ws.Range("D46").Value = "'0.0950"
$ws.Range("E46").Value = '  -0.86%  '

# Row 47
$ws.Range("D47").Value = "'0.0517"
$ws.Range("E47").Value = '  -1.73%  '

# Row 48
$ws.Range("D48").Value = "'18.40"

# Row 49
$ws.Range("D49").Value = "'1.969.44"
$ws.Range("E49").Value = '  -0.67%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = "'4.49"
$ws.Range("E50").Value = '  -3.88%  '

# Row 51
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = "'0.0221"
$ws.Range("E51").Value = '  -0.80%  '
